$d = $word.ActiveDocument

# 1. Heading: "You made our event a success! 🎉"
$d.Content.Find.Execute(
  "You made our event a success! 🎉", $true, $false, $false, $false, $false,
  $true, 1, $false,
  "Sự tham dự của bạn đã góp phần vào thành công của sự kiện chúng tôi! 🎉", 2)

# 2. "Hi " -> "Xin chào "
$d.Content.Find.Execute(
  "Hi ", $true, $false, $false, $false, $false,
  $true, 1, $false,
  "Xin chào ", 2)

# 3. "Thank you for attending " -> "Cảm ơn bạn đã tham dự sự kiện "
$d.Content.Find.Execute(
  "Thank you for attending ", $true, $false, $false, $false, $false,
  $true, 1, $false,
  "Cảm ơn bạn đã tham dự sự kiện ", 2)

# 4. " in " -> " tại "
$d.Content.Find.Execute(
  " in ", $true, $false, $false, $false, $false,
  $true, 1, $false,
  " tại ", 2)

# 5. ". We hope you had a great time, and it was a pleasure getting to know you!"
$d.Content.Find.Execute(
  ". We hope you had a great time, and it was a pleasure getting to know you!", $true, $false, $false, $false, $false,
  $true, 1, $false,
  ". Chúng tôi hy vọng bạn đã có một khoảng thời gian tuyệt vời. Rất vinh dự khi được làm quen với bạn!", 2)

# 6. "We hope the event inspired you as much as it did us, and let’s keep growing together!"
$d.Content.Find.Execute(
  "We hope the event inspired you as much as it did us, and let’s keep growing together!", $true, $false, $false, $false, $false,
  $true, 1, $false,
  "Chúng tôi mong rằng sự kiện đã có thể truyền cảm hứng cho bạn như nó đã làm với chúng tôi. Chúng tôi hy vọng bạn và chúng tôi sẽ cùng nhau phát triển hơn nữa trong tương lai!", 2)
